$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the hyperlink (email) cell style so we can re-apply it after the
# row shift below (inserting/re-adding hyperlinks can reset cell styling).
$hyperlinkStyle = $ws.Range("B1").Style

# Capture the existing hyperlink target addresses, in on-sheet order, before
# we touch anything - row Insert() does not relocate the Hyperlinks
# collection in this engine, so we rebuild it explicitly afterwards.
$linkAddrs = @()
foreach ($hh in $ws.Hyperlinks) {
    $linkAddrs += $hh.Address
}
$linkAddr1 = $linkAddrs[0]
$linkAddr2 = $linkAddrs[1]
$linkAddr3 = $linkAddrs[2]

# Insert a new row at the top to make room for the header row; this shifts
# the existing data (and cell styles/number formats) down one row.
$ws.Rows.Item(1).Insert()

# Set header values for the newly inserted row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "email"

# Rebuild the hyperlinks so they point at the data's new (shifted) rows.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), $linkAddr1)
$ws.Hyperlinks.Add($ws.Range("B3"), $linkAddr2)
$ws.Hyperlinks.Add($ws.Range("B4"), $linkAddr3)

# Adding hyperlinks can restyle the target cells - restore the original
# hyperlink look so only the header cell (B1) is left unstyled.
$ws.Range("B2").Style = $hyperlinkStyle
$ws.Range("B3").Style = $hyperlinkStyle
$ws.Range("B4").Style = $hyperlinkStyle

# Move the active selection down to the first empty row below the table,
# matching where Excel leaves the cursor after appending rows.
[void]$ws.Range("A5").Select()
